$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Item Code / Name / Purchase Rate / Sales Rate / Qty / Value between
# adjacent duplicate-name rows (columns B-G), leaving Sr.No (col A) untouched.

$ws.Range("B148").Value2 = 64196
$ws.Range("B149").Value2 = 65258
$ws.Range("F148").Value2 = 1
$ws.Range("F149").Value2 = 2
$ws.Range("G148").Value2 = 32143.58
$ws.Range("G149").Value2 = 64287.16
$ws.Range("B195").Value2 = 64350
$ws.Range("B196").Value2 = 57756
$ws.Range("E195").Value2 = 70.63
$ws.Range("E196").Value2 = 79.37
$ws.Range("F195").Value2 = 2
$ws.Range("F196").Value2 = -100
$ws.Range("G195").Value2 = 132.88
$ws.Range("G196").Value2 = -6644
$ws.Range("B233").Value2 = 57004
$ws.Range("B234").Value2 = 63255
$ws.Range("F233").Value2 = 5
$ws.Range("F234").Value2 = 82
$ws.Range("G233").Value2 = 410
$ws.Range("G234").Value2 = 6724
$ws.Range("B322").Value2 = 66188
$ws.Range("B323").Value2 = 48719
$ws.Range("C322").Value2 = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("C323").Value2 = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D322").Value2 = 315.8
$ws.Range("D323").Value2 = 295.75
$ws.Range("E322").Value2 = 377.31
$ws.Range("E323").Value2 = 353.35
$ws.Range("F322").Value2 = 36
$ws.Range("F323").Value2 = -82
$ws.Range("G322").Value2 = 11368.8
$ws.Range("G323").Value2 = -24251.5
$ws.Range("B367").Value2 = 66194
$ws.Range("B368").Value2 = 64983
$ws.Range("C367").Value2 = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("C368").Value2 = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F367").Value2 = 35
$ws.Range("F368").Value2 = 6
$ws.Range("G367").Value2 = 2998.8
$ws.Range("G368").Value2 = 514.08
$ws.Range("B369").Value2 = 64985
$ws.Range("B370").Value2 = 66196
$ws.Range("C369").Value2 = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("C370").Value2 = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F369").Value2 = 13
$ws.Range("F370").Value2 = 28
$ws.Range("G369").Value2 = 1140.1
$ws.Range("G370").Value2 = 2455.6
$ws.Range("B375").Value2 = 63565
$ws.Range("B376").Value2 = 61610
$ws.Range("E375").Value2 = 109.19
$ws.Range("E376").Value2 = 122.71
$ws.Range("F375").Value2 = 60
$ws.Range("F376").Value2 = -58
$ws.Range("G375").Value2 = 6162.6
$ws.Range("G376").Value2 = -5957.18
$ws.Range("B397").Value2 = 63560
$ws.Range("B398").Value2 = 60325
$ws.Range("E397").Value2 = 134.87
$ws.Range("E398").Value2 = 151.57
$ws.Range("F397").Value2 = 1
$ws.Range("F398").Value2 = -102
$ws.Range("G397").Value2 = 126.86
$ws.Range("G398").Value2 = -12939.72
$ws.Range("B548").Value2 = 53602
$ws.Range("B549").Value2 = 65068
$ws.Range("E548").Value2 = 15.69
$ws.Range("E549").Value2 = 13.97
$ws.Range("F548").Value2 = -231
$ws.Range("F549").Value2 = 63
$ws.Range("G548").Value2 = -3037.65
$ws.Range("G549").Value2 = 828.45
$ws.Range("B556").Value2 = 45706
$ws.Range("B557").Value2 = 64922
$ws.Range("E556").Value2 = 23.58
$ws.Range("E557").Value2 = 20.98
$ws.Range("F556").Value2 = -202
$ws.Range("F557").Value2 = 67
$ws.Range("G556").Value2 = -3985.46
$ws.Range("G557").Value2 = 1321.91
$ws.Range("B559").Value2 = 64927
$ws.Range("B560").Value2 = 45718
$ws.Range("E559").Value2 = 17.26
$ws.Range("E560").Value2 = 19.38
$ws.Range("F559").Value2 = 106
$ws.Range("F560").Value2 = -294
$ws.Range("G559").Value2 = 1719.32
$ws.Range("G560").Value2 = -4768.68
$ws.Range("B564").Value2 = 64925
$ws.Range("B565").Value2 = 45709
$ws.Range("E564").Value2 = 13.97
$ws.Range("E565").Value2 = 15.69
$ws.Range("F564").Value2 = 111
$ws.Range("F565").Value2 = -300
$ws.Range("G564").Value2 = 1459.65
$ws.Range("G565").Value2 = -3945
$ws.Range("B566").Value2 = 45702
$ws.Range("B567").Value2 = 64919
$ws.Range("E566").Value2 = 31.43
$ws.Range("E567").Value2 = 27.97
$ws.Range("F566").Value2 = -215
$ws.Range("F567").Value2 = 61
$ws.Range("G566").Value2 = -5654.5
$ws.Range("G567").Value2 = 1604.3
$ws.Range("B569").Value2 = 53595
$ws.Range("B570").Value2 = 65067
$ws.Range("E569").Value2 = 17.61
$ws.Range("E570").Value2 = 15.65
$ws.Range("F569").Value2 = -335
$ws.Range("F570").Value2 = 126
$ws.Range("G569").Value2 = -4934.55
$ws.Range("G570").Value2 = 1855.98
$ws.Range("B640").Value2 = 64810
$ws.Range("B641").Value2 = 53319
$ws.Range("E640").Value2 = 291.22
$ws.Range("E641").Value2 = 310.64
$ws.Range("F640").Value2 = 2
$ws.Range("F641").Value2 = -6
$ws.Range("G640").Value2 = 547.84
$ws.Range("G641").Value2 = -1643.52
$ws.Range("B659").Value2 = 64833
$ws.Range("B660").Value2 = 60025
$ws.Range("E659").Value2 = 34.9
$ws.Range("E660").Value2 = 37.22
$ws.Range("F659").Value2 = 88
$ws.Range("F660").Value2 = -98
$ws.Range("G659").Value2 = 2889.04
$ws.Range("G660").Value2 = -3217.34
